# Add diary entries for 2/13, 2/16, 2/18, 2/19 and 2/20 into rows 19-23.
# (the 7 leftover placeholder rows at 19-25 shrink down to 2 placeholder
# rows, now at 24-25, after the 5 new entries are inserted)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 19: 2/13/2020 ----
$ws.Range("A19").Value = 43874
$ws.Range("A19").NumberFormat = "MM/DD/YY"
$ws.Range("B19").Value = "5PM-8PM"
$ws.Range("C19").Value = "Class"
$ws.Range("D19").Value = "Attend lecture, take midterm"
$ws.Range("E19").Value = "All goals"
$ws.Range("F19").Value = "Explaining things is hard."
$ws.Range("G19").Value = "Stressed"
$ws.Rows.Item(19).RowHeight = 15.65

# ---- Row 20: 2/16/2020 ----
$ws.Range("A20").Value = 43877
$ws.Range("A20").NumberFormat = "MM/DD/YY"
$ws.Range("B20").Value = "3PM-7PM"
$ws.Range("C20").Value = "Team"
$ws.Range("D20").Value = "Finish group project assignment, discuss assignment 2 redo"
$ws.Range("E20").Value = "All goals"
$ws.Range("F20").Value = "Much of the information on our open-source project’s goals, stakeholders, etc. was available on their website. There was also a whole community of players and mailing lists that we were unaware of before."
$ws.Range("G20").Value = "Happy & pleasantly surprised"
$ws.Rows.Item(20).RowHeight = 86.55

# ---- Row 21: 2/18/2020 ----
$ws.Range("A21").Value = 43879
$ws.Range("A21").NumberFormat = "MM/DD/YY"
$ws.Range("B21").Value = "4PM-4:30PM"
$ws.Range("C21").Value = "Team, Kaj"
$ws.Range("D21").Value = "Go over feedback for assignment 2"
$ws.Range("E21").Value = "All goals"
$ws.Range("F21").Value = "It was helpful to be reminded of the importance of explaining diagrams and images, which was a blind spot for me. I also appreciated knowing more about the level of abstraction expected for this kind of documentation."
$ws.Range("G21").Value = "Thanks Kaj!"
$ws.Rows.Item(21).RowHeight = 100.7

# ---- Row 22: 2/19/2020 ----
$ws.Range("A22").Value = 43880
$ws.Range("A22").NumberFormat = "MM/DD/YY"
$ws.Range("B22").Value = "2:30PM-5PM"
$ws.Range("C22").Value = "Team"
$ws.Range("D22").Value = "Redo assignment 2"
$ws.Range("E22").Value = "Almost done, just missing a few more details"
$ws.Range("F22").Value = "It is difficult to determine what level of detail to go into when describing something at a high level. We tried to stick to the essentials and omit minor details that would cloud the big-picture view of the features"
$ws.Range("G22").Value = "Tired"
$ws.Rows.Item(22).RowHeight = 86.55

# ---- Row 23: 2/20/2020 ----
$ws.Range("A23").Value = 43881
$ws.Range("A23").NumberFormat = "MM/DD/YY"
$ws.Range("B23").Value = "5PM-8PM"
$ws.Range("C23").Value = "Class"
$ws.Range("D23").Value = "Attend lecture"
$ws.Range("E23").Value = "Learned about reading code in terms of architectural styles and social context"
$ws.Range("F23").Value = "Determining what would be useful in terms of illustrating the architecture of a program is pretty hard. It is also pretty easy to fall into the mistake of assuming an architectural style is being used when that’s not the case. Seeing some of the ways to find metrics for determining the state of an open source project was helpful."
$ws.Range("G23").Value = "Positive"
$ws.Rows.Item(23).RowHeight = 129.1

# The Achievements cell of the last new entry gets its own (visually
# identical) italic style slot, same as in the source workbook.
$f23 = $ws.Range("E23").Font
$f23.Italic = $true

# Update selection / scroll position to match the edited area.
$ws.Range("A24").Select()
